$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 26 / 27 swap (LidoDAOToken <-> EthereumClassic)
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

# Row 45 / 46 swap (EnergySwap <-> RenderToken)
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

# Column D (Price) updates - forced to text to preserve literal formatting
Set-TextValue $ws.Range("D2") "26.875.20"
Set-TextValue $ws.Range("D3") "1.807.35"
Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("D5") "310.89"
Set-TextValue $ws.Range("D7") "0.4463"
Set-TextValue $ws.Range("D8") "0.3674"
Set-TextValue $ws.Range("D9") "0.07410"
Set-TextValue $ws.Range("D10") "0.8560"
Set-TextValue $ws.Range("D11") "20.69"
Set-TextValue $ws.Range("D12") "1.806.81"
Set-TextValue $ws.Range("D13") "6.608"
Set-TextValue $ws.Range("D14") "92.75"
Set-TextValue $ws.Range("D15") "5.311"
Set-TextValue $ws.Range("D16") "0.07071"
Set-TextValue $ws.Range("D18") "0.000008749"
Set-TextValue $ws.Range("D20") "14.90"
Set-TextValue $ws.Range("D21") "26.895.79"
Set-TextValue $ws.Range("D22") "5.157"
Set-TextValue $ws.Range("D23") "10.84"
Set-TextValue $ws.Range("D24") "1.991"
Set-TextValue $ws.Range("D25") "151.97"
Set-TextValue $ws.Range("D26") "18.50"
Set-TextValue $ws.Range("D27") "2.184"
Set-TextValue $ws.Range("D28") "5.211"
Set-TextValue $ws.Range("D29") "116.53"
Set-TextValue $ws.Range("D30") "0.08839"
Set-TextValue $ws.Range("D31") "0.7544"
Set-TextValue $ws.Range("D32") "1.176"
Set-TextValue $ws.Range("D33") "2.925"
Set-TextValue $ws.Range("D34") "4.463"
Set-TextValue $ws.Range("D36") "1.090"
Set-TextValue $ws.Range("D37") "0.01972"
Set-TextValue $ws.Range("D38") "0.05205"
Set-TextValue $ws.Range("D40") "2.866"
Set-TextValue $ws.Range("D41") "7.010"
Set-TextValue $ws.Range("D43") "0.5182"
Set-TextValue $ws.Range("D44") "8.444"
Set-TextValue $ws.Range("D45") "1.990"
Set-TextValue $ws.Range("D46") "10.52"
Set-TextValue $ws.Range("D47") "105.50"
Set-TextValue $ws.Range("D49") "0.9994"
Set-TextValue $ws.Range("D50") "0.06331"
Set-TextValue $ws.Range("D51") "0.9215"

# Column E (Volume 1h) updates - forced to text to preserve literal formatting/spacing
Set-TextValue $ws.Range("E2") "  -1.22%  "
Set-TextValue $ws.Range("E3") "  -0.99%  "
Set-TextValue $ws.Range("E4") "  -0.48%  "
Set-TextValue $ws.Range("E5") "  -0.80%  "
Set-TextValue $ws.Range("E6") "  -0.43%  "
Set-TextValue $ws.Range("E7") "  +5.08%  "
Set-TextValue $ws.Range("E8") "  -0.79%  "
Set-TextValue $ws.Range("E9") "  +2.17%  "
Set-TextValue $ws.Range("E10") "  -0.84%  "
Set-TextValue $ws.Range("E11") "  -1.82%  "
Set-TextValue $ws.Range("E12") "  -1.02%  "
Set-TextValue $ws.Range("E13") "  -1.77%  "
Set-TextValue $ws.Range("E14") "  +3.40%  "
Set-TextValue $ws.Range("E15") "  -0.23%  "
Set-TextValue $ws.Range("E16") "  -0.29%  "
Set-TextValue $ws.Range("E17") "  -0.50%  "
Set-TextValue $ws.Range("E18") "  -1.24%  "
Set-TextValue $ws.Range("E19") "  -0.38%  "
Set-TextValue $ws.Range("E20") "  -1.34%  "
Set-TextValue $ws.Range("E21") "  -1.38%  "
Set-TextValue $ws.Range("E22") "  +0.32%  "
Set-TextValue $ws.Range("E23") "  -0.65%  "
Set-TextValue $ws.Range("E24") "  +0.15%  "
Set-TextValue $ws.Range("E26") "  +0.46%  "
Set-TextValue $ws.Range("E27") "  -0.94%  "
Set-TextValue $ws.Range("E28") "  -0.63%  "
Set-TextValue $ws.Range("E29") "  +0.12%  "
Set-TextValue $ws.Range("E30") "  +0.06%  "
Set-TextValue $ws.Range("E31") "  -0.26%  "
Set-TextValue $ws.Range("E32") "  -1.66%  "
Set-TextValue $ws.Range("E33") "  +6.50%  "
Set-TextValue $ws.Range("E34") "  +0.17%  "
Set-TextValue $ws.Range("E35") "  -0.51%  "
Set-TextValue $ws.Range("E36") "  -2.32%  "
Set-TextValue $ws.Range("E37") "  -0.05%  "
Set-TextValue $ws.Range("E38") "  -1.11%  "
Set-TextValue $ws.Range("E39") "  +5.44%  "
Set-TextValue $ws.Range("E40") "  -0.23%  "
Set-TextValue $ws.Range("E41") "  -4.73%  "
Set-TextValue $ws.Range("E42") "  -0.38%  "
Set-TextValue $ws.Range("E43") "  +9.06%  "
Set-TextValue $ws.Range("E44") "  -3.53%  "
Set-TextValue $ws.Range("E45") "  +6.93%  "
Set-TextValue $ws.Range("E46") "  -1.22%  "
Set-TextValue $ws.Range("E47") "  -1.83%  "
Set-TextValue $ws.Range("E48") "  -0.04%  "
Set-TextValue $ws.Range("E49") "  -0.52%  "
Set-TextValue $ws.Range("E50") "  -0.91%  "
Set-TextValue $ws.Range("E51") "  +0.78%  "
